$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated Price (column D) and Volume(1h) (column E) values.
# Leading apostrophes force text interpretation so values are stored
# the same way as the original inline-string cells (e.g. percentages
# and numbers keep their exact original textual formatting).
$ws.Range("D2").Value = "'276.87"
$ws.Range("E2").Value = "'0.89%"
$ws.Range("D3").Value = "'27.11"
$ws.Range("E3").Value = "'1.24%"
$ws.Range("D4").Value = "'4.855"
$ws.Range("E4").Value = "'0.12%"
$ws.Range("D5").Value = "'0.06406"
$ws.Range("E5").Value = "'1.32%"
$ws.Range("D6").Value = "'6.934"
$ws.Range("E6").Value = "'0.72%"
$ws.Range("D7").Value = "'1.196"
$ws.Range("E7").Value = "'-6.54%"
$ws.Range("D8").Value = "'0.8766"
$ws.Range("E8").Value = "'0.68%"
$ws.Range("D9").Value = "'0.1544"
$ws.Range("E9").Value = "'5.69%"
$ws.Range("D10").Value = "'0.05130"
$ws.Range("E10").Value = "'2.74%"
$ws.Range("D11").Value = "'0.07486"
$ws.Range("E11").Value = "'1.13%"
$ws.Range("D12").Value = "'0.02965"
$ws.Range("E12").Value = "'0.80%"
$ws.Range("D13").Value = "'0.08978"
$ws.Range("E13").Value = "'-0.54%"
$ws.Range("D14").Value = "'0.001566"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("D15").Value = "'0.0006394"
$ws.Range("E15").Value = "'1.62%"
$ws.Range("D16").Value = "'0.006086"
$ws.Range("E16").Value = "'1.09%"
$ws.Range("D17").Value = "'3.478"
$ws.Range("E17").Value = "'0.94%"
$ws.Range("D18").Value = "'3.309"
$ws.Range("E18").Value = "'-0.43%"
$ws.Range("E19").Value = "'-0.39%"
$ws.Range("E20").Value = "'0.29%"
$ws.Range("E21").Value = "'0.94%"
$ws.Range("D22").Value = "'3.903"
$ws.Range("E22").Value = "'-0.32%"
$ws.Range("D23").Value = "'0.04426"
$ws.Range("E23").Value = "'1.73%"
$ws.Range("D25").Value = "'0.001178"
$ws.Range("E25").Value = "'-0.13%"
$ws.Range("D26").Value = "'0.003860"
$ws.Range("E26").Value = "'-9.26%"
$ws.Range("E28").Value = "'15.09%"
$ws.Range("D40").Value = "'0.04160"
$ws.Range("E40").Value = "'2.84%"
$ws.Range("E41").Value = "'1.77%"
$ws.Range("E42").Value = "'0.74%"
$ws.Range("D43").Value = "'0.001950"
$ws.Range("E43").Value = "'-7.15%"
$ws.Range("D44").Value = "'0.01191"
$ws.Range("E44").Value = "'11.33%"
$ws.Range("D45").Value = "'0.00005305"
$ws.Range("E45").Value = "'0.01%"
$ws.Range("E46").Value = "'13.51%"
$ws.Range("D47").Value = "'0.01853"
$ws.Range("E47").Value = "'-7.39%"
